$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Randomize the prime/target image pairings in columns A and B (rows 2-61)
$values = @(
    @(2, "Primes/C-1.jpg", "Targets/T (1).BMP"),
    @(3, "Primes/C-2.jpg", "Targets/T (2).BMP"),
    @(4, "Primes/S-7.jpg", "Targets/T (37).BMP"),
    @(5, "Primes/C-6.jpg", "Targets/T (26).BMP"),
    @(6, "Primes/S-5.jpg", "Targets/T (45).BMP"),
    @(7, "Primes/S-10.jpg", "Targets/T (40).BMP"),
    @(8, "Primes/C-7.jpg", "Targets/T (7).BMP"),
    @(9, "Primes/S-5.jpg", "Targets/T (55).BMP"),
    @(10, "Primes/S-6.jpg", "Targets/T (56).BMP"),
    @(11, "Primes/C-10.jpg", "Targets/T (10).BMP"),
    @(12, "Primes/S-4.jpg", "Targets/T (34).BMP"),
    @(13, "Primes/S-5.jpg", "Targets/T (35).BMP"),
    @(14, "Primes/C-3.jpg", "Targets/T (13).BMP"),
    @(15, "Primes/C-4.jpg", "Targets/T (14).BMP"),
    @(16, "Primes/S-8.jpg", "Targets/T (58).BMP"),
    @(17, "Primes/C-8.jpg", "Targets/T (18).BMP"),
    @(18, "Primes/S-10.jpg", "Targets/T (60).BMP"),
    @(19, "Primes/S-9.jpg", "Targets/T (59).BMP"),
    @(20, "Primes/C-9.jpg", "Targets/T (19).BMP"),
    @(21, "Primes/C-10.jpg", "Targets/T (20).BMP"),
    @(22, "Primes/C-1.jpg", "Targets/T (21).BMP"),
    @(23, "Primes/S-7.jpg", "Targets/T (47).BMP"),
    @(24, "Primes/S-3.jpg", "Targets/T (33).BMP"),
    @(25, "Primes/S-9.jpg", "Targets/T (49).BMP"),
    @(26, "Primes/C-1.jpg", "Targets/T (11).BMP"),
    @(27, "Primes/S-10.jpg", "Targets/T (50).BMP"),
    @(28, "Primes/S-2.jpg", "Targets/T (52).BMP"),
    @(29, "Primes/S-3.jpg", "Targets/T (53).BMP"),
    @(30, "Primes/C-9.jpg", "Targets/T (29).BMP"),
    @(31, "Primes/C-10.jpg", "Targets/T (30).BMP"),
    @(32, "Primes/S-1.jpg", "Targets/T (31).BMP"),
    @(33, "Primes/S-2.jpg", "Targets/T (32).BMP"),
    @(34, "Primes/S-6.jpg", "Targets/T (36).BMP"),
    @(35, "Primes/C-2.jpg", "Targets/T (12).BMP"),
    @(36, "Primes/S-8.jpg", "Targets/T (48).BMP"),
    @(37, "Primes/S-1.jpg", "Targets/T (51).BMP"),
    @(38, "Primes/C-3.jpg", "Targets/T (3).BMP"),
    @(39, "Primes/C-6.jpg", "Targets/T (6).BMP"),
    @(40, "Primes/C-5.jpg", "Targets/T (5).BMP"),
    @(41, "Primes/S-1.jpg", "Targets/T (41).BMP"),
    @(42, "Primes/C-3.jpg", "Targets/T (23).BMP"),
    @(43, "Primes/S-2.jpg", "Targets/T (42).BMP"),
    @(44, "Primes/C-5.jpg", "Targets/T (25).BMP"),
    @(45, "Primes/S-8.jpg", "Targets/T (38).BMP"),
    @(46, "Primes/S-9.jpg", "Targets/T (39).BMP"),
    @(47, "Primes/S-6.jpg", "Targets/T (46).BMP"),
    @(48, "Primes/C-2.jpg", "Targets/T (22).BMP"),
    @(49, "Primes/C-4.jpg", "Targets/T (4).BMP"),
    @(50, "Primes/C-4.jpg", "Targets/T (24).BMP"),
    @(51, "Primes/S-3.jpg", "Targets/T (43).BMP"),
    @(52, "Primes/S-4.jpg", "Targets/T (44).BMP"),
    @(53, "Primes/C-7.jpg", "Targets/T (27).BMP"),
    @(54, "Primes/C-8.jpg", "Targets/T (28).BMP"),
    @(55, "Primes/S-4.jpg", "Targets/T (54).BMP"),
    @(56, "Primes/C-8.jpg", "Targets/T (8).BMP"),
    @(57, "Primes/C-9.jpg", "Targets/T (9).BMP"),
    @(58, "Primes/S-7.jpg", "Targets/T (57).BMP"),
    @(59, "Primes/C-5.jpg", "Targets/T (15).BMP"),
    @(60, "Primes/C-6.jpg", "Targets/T (16).BMP"),
    @(61, "Primes/C-7.jpg", "Targets/T (17).BMP")
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Rows that previously used the special "Neutral"/accent styles now use the
# same plain centered style as the rest of column A (copy format from A4).
$ws.Range("A4").Copy() | Out-Null
$styleFixRows = @(2,3,12,13,22,23,32,42,52)
foreach ($r in $styleFixRows) {
    $ws.Range("A$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# The "Neutral" cell style is no longer used anywhere in the sheet; remove it
# from the workbook's style list.
$wb.Styles("Neutral").Delete()
